$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)
$sh = $s.Shapes.Item(2)

# --- Paragraph: "Page type counts:" at the same indent level (lvl=2 -> IndentLevel 3)
# as the "enwiki-latest-pages-articles10.xml-p000925001p001325000" paragraph (paragraph 3).
$tr = $sh.TextFrame.TextRange
$srcPara = $tr.Paragraphs(3)
$srcPara.InsertAfter("`rPage type counts:")

# Re-fetch the freshly created paragraph (now paragraph 4) and split it into three runs
# "Page " / "type counts" / ":" by rewriting sub-ranges (matches the run boundaries
# produced when PowerPoint commits autocorrect/spellcheck passes while typing).
$tr = $sh.TextFrame.TextRange
$newPara = $tr.Paragraphs(4)
$colonRange = $newPara.Characters(17, 1)
$colonRange.Text = ":"
$typeCountsRange = $newPara.Characters(6, 11)
$typeCountsRange.Text = "type counts"

# --- Paragraph: "ARTICLE 178045" one level deeper (lvl=3 -> IndentLevel 4)
$tr = $sh.TextFrame.TextRange
$newPara = $tr.Paragraphs(4)
$newPara.InsertAfter("`rARTICLE 178045")
$tr = $sh.TextFrame.TextRange
$articlePara = $tr.Paragraphs(5)
$articlePara.IndentLevel = 4

# --- Paragraph: "DISAMBIGUATION 5454" at the same deeper level
$tr = $sh.TextFrame.TextRange
$articlePara = $tr.Paragraphs(5)
$articlePara.InsertAfter("`rDISAMBIGUATION 5454")
$tr = $sh.TextFrame.TextRange
$disambigPara = $tr.Paragraphs(6)
$disambigPara.IndentLevel = 4

# --- Paragraph: "NON_ARTICLE 31963" at the same deeper level
$tr = $sh.TextFrame.TextRange
$disambigPara = $tr.Paragraphs(6)
$disambigPara.InsertAfter("`rNON_ARTICLE 31963")
$tr = $sh.TextFrame.TextRange
$nonArticlePara = $tr.Paragraphs(7)
$nonArticlePara.IndentLevel = 4
